$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Quantity (column B) and Order Number (column D) for rows 2 and 3.
# Both columns hold text values in the source data (e.g. leading zeros in
# Order Number), so force the cell format to Text before assigning values
# to keep them stored as strings rather than numbers.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "9"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "000001305"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "9"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "000001305"
